$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the birthday values (dates stored as day-count serials)
$ws1.Range("B4").Value = (Get-Date -Year 1991 -Month 3 -Day 27).Date
$ws1.Range("B6").Value = (Get-Date -Year 1989 -Month 3 -Day 26).Date

# Update the active selection on Sheet1
$ws1.Activate()
$ws1.Range("E11").Select()
